$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets: 1 = Overview, 2 = zh-cn, 3 = de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0e7fb31fe155c2459a1939fca5e169d000ed1f5/e2e/62727eb8-e36d-4f5d-a19d-6325ddf64c87.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0e7fb31fe155c2459a1939fca5e169d000ed1f5/e2e/bbc82d48-1963-4fa7-ae16-6fd5d8f2f120.md"
$mdName1 = "62727eb8-e36d-4f5d-a19d-6325ddf64c87.md"
$mdName2 = "bbc82d48-1963-4fa7-ae16-6fd5d8f2f120.md"

$hyperlinkColor = 15570276   # OLE BGR for RGB(0x64,0x95,0xED) -- the same blue used by the workbook's existing HyperLink style

function Apply-HandbackRow($ws, $mdUrl, $mdName, $handbackFileValue, $row) {
    $iCell = $ws.Range("I" + $row)
    $ws.Hyperlinks.Add($iCell, $mdUrl, "", "", $mdName) | Out-Null
    $iCell.Font.Underline = $true
    $iCell.Font.Color = $hyperlinkColor

    $ws.Range("J" + $row).Value = $handbackFileValue
}

# ---------------------------------------------------------------------------
# Status column (C) -> "Handed back: in sync with en-US" for both data rows,
# on both the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in "Latest Target File" (I), "Latest Handback File" (J)
# and "Latest Handback DateTime" (K) for both rows.
# ---------------------------------------------------------------------------
Apply-HandbackRow $wsZhCn $mdUrl1 $mdName1 "62727eb8-e36d-4f5d-a19d-6325ddf64c87.6356c7d22766ca2be6fd2c23deda3e0edfce368f.zh-cn.xlf" 2
Apply-HandbackRow $wsZhCn $mdUrl2 $mdName2 "bbc82d48-1963-4fa7-ae16-6fd5d8f2f120.faeb4e0617927451cca940a621665e48fbdc47b0.zh-cn.xlf" 3

$wsZhCn.Range("K2").Value = "2016-09-01 22:31:40"
$wsZhCn.Range("K3").Value = "2016-09-01 22:31:40"

# ---------------------------------------------------------------------------
# de-de sheet: same treatment, different handback datetime stamp.
# ---------------------------------------------------------------------------
Apply-HandbackRow $wsDeDe $mdUrl1 $mdName1 "62727eb8-e36d-4f5d-a19d-6325ddf64c87.6356c7d22766ca2be6fd2c23deda3e0edfce368f.de-de.xlf" 2
Apply-HandbackRow $wsDeDe $mdUrl2 $mdName2 "bbc82d48-1963-4fa7-ae16-6fd5d8f2f120.faeb4e0617927451cca940a621665e48fbdc47b0.de-de.xlf" 3

$wsDeDe.Range("K2").Value = "2016-09-01 22:31:48"
$wsDeDe.Range("K3").Value = "2016-09-01 22:31:48"

# ---------------------------------------------------------------------------
# Column widths — widened to fit the newly populated columns.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.1    # E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = 29.1    # F: de-de

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1        # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth = 39.14       # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.14      # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1        # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth = 39.14       # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.14      # J: Latest Handback File
